# Auto-generated edit script: applies per-cell numeric updates to match the target diff.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across several
# worksheets, matching a refreshed market-data pull from the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value2 = 156416.08
$ws.Range("J75").Value2 = 120744.75
$ws.Range("L75").Value2 = 120744.75
$ws.Range("N75").Value2 = -122616.75
$ws.Range("H78").Value2 = 156416.08
$ws.Range("J78").Value2 = 120744.75
$ws.Range("L78").Value2 = 362234.25
$ws.Range("N78").Value2 = -371594.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 3286754.2
$ws.Range("I2").Value2 = 4382101
$ws.Range("J2").Value2 = 713.5714
$ws.Range("K2").Value2 = 4382101
$ws.Range("L2").Value2 = 713.5714
$ws.Range("M2").Value2 = -4381988
$ws.Range("N2").Value2 = -939.5714
$ws.Range("H45").Value2 = 1925.25
$ws.Range("I45").Value2 = 1892.8
$ws.Range("K45").Value2 = 1892.8
$ws.Range("M45").Value2 = -1515.8
$ws.Range("H74").Value2 = 457285.53
$ws.Range("I74").Value2 = 834562.9399999999
$ws.Range("K74").Value2 = 834562.9399999999
$ws.Range("M74").Value2 = -833688.9399999999
$ws.Range("H77").Value2 = 457285.53
$ws.Range("I77").Value2 = 834562.9399999999
$ws.Range("K77").Value2 = 4172814.7
$ws.Range("M77").Value2 = -4168446.7
$ws.Range("H110").Value2 = 35715710
$ws.Range("I110").Value2 = 62500270
$ws.Range("K110").Value2 = 62500270
$ws.Range("M110").Value2 = -62498225
$ws.Range("H116").Value2 = 3286754.2
$ws.Range("I116").Value2 = 4382101
$ws.Range("J116").Value2 = 713.5714
$ws.Range("K116").Value2 = 4382101
$ws.Range("L116").Value2 = 713.5714
$ws.Range("M116").Value2 = -4379807
$ws.Range("N116").Value2 = -5301.5714
$ws.Range("H132").Value2 = 5720704.5
$ws.Range("I132").Value2 = 12504697
$ws.Range("K132").Value2 = 37514091
$ws.Range("M132").Value2 = -37511561

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 3286754.2
$ws.Range("I3").Value2 = 4382101
$ws.Range("J3").Value2 = 713.5714
$ws.Range("K3").Value2 = 4382101
$ws.Range("L3").Value2 = 713.5714
$ws.Range("M3").Value2 = -4381987
$ws.Range("N3").Value2 = -941.5714
$ws.Range("H134").Value2 = 7037.52
$ws.Range("I134").Value2 = 5392.2
$ws.Range("J134").Value2 = 9505.5
$ws.Range("K134").Value2 = 16176.6
$ws.Range("L134").Value2 = 28516.5
$ws.Range("M134").Value2 = -13641.6
$ws.Range("N134").Value2 = -33586.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 5872.1055
$ws.Range("I31").Value2 = 3478.5
$ws.Range("K31").Value2 = 3478.5
$ws.Range("M31").Value2 = -3183.5
$ws.Range("H34").Value2 = 5872.1055
$ws.Range("I34").Value2 = 3478.5
$ws.Range("K34").Value2 = 3478.5
$ws.Range("M34").Value2 = -3276.5
$ws.Range("H88").Value2 = 20667
$ws.Range("J88").Value2 = 25334
$ws.Range("L88").Value2 = 25334
$ws.Range("N88").Value2 = -26146
$ws.Range("H91").Value2 = 20667
$ws.Range("J91").Value2 = 25334
$ws.Range("L91").Value2 = 25334
$ws.Range("N91").Value2 = -28142
$ws.Range("H99").Value2 = 4018.1333
$ws.Range("I99").Value2 = 3615.375
$ws.Range("K99").Value2 = 3615.375
$ws.Range("M99").Value2 = -2117.375
$ws.Range("H126").Value2 = 4018.1333
$ws.Range("I126").Value2 = 3615.375
$ws.Range("K126").Value2 = 10846.125
$ws.Range("M126").Value2 = -8376.125
$ws.Range("H132").Value2 = 41803.438
$ws.Range("I132").Value2 = 4420.857
$ws.Range("J132").Value2 = 70878.78
$ws.Range("K132").Value2 = 13262.571
$ws.Range("L132").Value2 = 212636.34
$ws.Range("M132").Value2 = -10732.571
$ws.Range("N132").Value2 = -217696.34

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value2 = 332.66666
$ws.Range("I25").Value2 = 332.66666
$ws.Range("K25").Value2 = 997.9999799999999
$ws.Range("M25").Value2 = -828.9999799999999
$ws.Range("H29").Value2 = 0
$ws.Range("J29").Value2 = 0
$ws.Range("L29").Value2 = 0
$ws.Range("N29").ClearContents()
$ws.Range("H30").Value2 = 332.66666
$ws.Range("I30").Value2 = 332.66666
$ws.Range("K30").Value2 = 997.9999799999999
$ws.Range("M30").Value2 = -895.9999799999999
$ws.Range("H31").Value2 = 0
$ws.Range("I31").Value2 = 0
$ws.Range("K31").Value2 = 0
$ws.Range("M31").ClearContents()
$ws.Range("H35").Value2 = 1199.7142
$ws.Range("J35").Value2 = 1774.5
$ws.Range("L35").Value2 = 5323.5
$ws.Range("N35").Value2 = -5899.5
$ws.Range("H80").Value2 = 3273.625
$ws.Range("J80").Value2 = 3470.7144
$ws.Range("L80").Value2 = 10412.1432
$ws.Range("N80").Value2 = -12284.1432
$ws.Range("H83").Value2 = 3273.625
$ws.Range("J83").Value2 = 3470.7144
$ws.Range("L83").Value2 = 31236.4296
$ws.Range("N83").Value2 = -40596.4296
$ws.Range("H86").Value2 = 279.77777
$ws.Range("I86").Value2 = 498
$ws.Range("J86").Value2 = 170.66667
$ws.Range("K86").Value2 = 1494
$ws.Range("L86").Value2 = 512.00001
$ws.Range("M86").Value2 = -308
$ws.Range("N86").Value2 = -2884.00001
$ws.Range("H89").Value2 = 279.77777
$ws.Range("I89").Value2 = 498
$ws.Range("J89").Value2 = 170.66667
$ws.Range("K89").Value2 = 4482
$ws.Range("L89").Value2 = 1536.00003
$ws.Range("M89").Value2 = 1446
$ws.Range("N89").Value2 = -13392.00003
$ws.Range("H113").Value2 = 2608.3572
$ws.Range("I113").Value2 = 2974.8
$ws.Range("J113").Value2 = 2404.7778
$ws.Range("K113").Value2 = 8924.400000000001
$ws.Range("L113").Value2 = 7214.3334
$ws.Range("M113").Value2 = -6754.400000000001
$ws.Range("N113").Value2 = -11554.3334
$ws.Range("H140").Value2 = 1704
$ws.Range("I140").Value2 = 1165
$ws.Range("J140").Value2 = 1906.125
$ws.Range("K140").Value2 = 3495
$ws.Range("L140").Value2 = 5718.375
$ws.Range("M140").Value2 = 1685
$ws.Range("N140").Value2 = -16078.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value2 = 1233.8572
$ws.Range("I113").Value2 = 996.75
$ws.Range("K113").Value2 = 996.75
$ws.Range("M113").Value2 = 1173.25
$ws.Range("H126").Value2 = 3125.0667
$ws.Range("I126").Value2 = 2510.9092
$ws.Range("K126").Value2 = 7532.7276
$ws.Range("M126").Value2 = -5062.7276
$ws.Range("H132").Value2 = 6225.788
$ws.Range("I132").Value2 = 4409.5
$ws.Range("K132").Value2 = 13228.5
$ws.Range("M132").Value2 = -10698.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value2 = 125012000
$ws.Range("I40").Value2 = 125012000
$ws.Range("K40").Value2 = 125012000
$ws.Range("M40").Value2 = -125011864
$ws.Range("H82").Value2 = 1996
$ws.Range("I82").Value2 = 1997.4286
$ws.Range("J82").Value2 = 1993.5
$ws.Range("K82").Value2 = 1997.4286
$ws.Range("L82").Value2 = 1993.5
$ws.Range("M82").Value2 = -1636.4286
$ws.Range("N82").Value2 = -2715.5
$ws.Range("H85").Value2 = 1996
$ws.Range("I85").Value2 = 1997.4286
$ws.Range("J85").Value2 = 1993.5
$ws.Range("K85").Value2 = 1997.4286
$ws.Range("L85").Value2 = 1993.5
$ws.Range("M85").Value2 = -749.4286
$ws.Range("N85").Value2 = -4489.5
$ws.Range("H122").Value2 = 3592.5
$ws.Range("I122").Value2 = 3290.1667
$ws.Range("K122").Value2 = 9870.500100000001
$ws.Range("M122").Value2 = -7420.500100000001
$ws.Range("H132").Value2 = 5963.773
$ws.Range("I132").Value2 = 4517
$ws.Range("K132").Value2 = 13551
$ws.Range("M132").Value2 = -11021
$ws.Range("H136").Value2 = 7415593
$ws.Range("I136").Value2 = 9530662
$ws.Range("K136").Value2 = 28591986
$ws.Range("M136").Value2 = -28589436

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value2 = 5060.4146
$ws.Range("I132").Value2 = 4630.1177
$ws.Range("J132").Value2 = 7150.4287
$ws.Range("K132").Value2 = 13890.3531
$ws.Range("L132").Value2 = 21451.2861
$ws.Range("M132").Value2 = -11360.3531
$ws.Range("N132").Value2 = -26511.2861
$ws.Range("H136").Value2 = 3704.261
$ws.Range("I136").Value2 = 1915.9333
$ws.Range("J136").Value2 = 7057.375
$ws.Range("K136").Value2 = 5747.7999
$ws.Range("L136").Value2 = 21172.125
$ws.Range("M136").Value2 = -3197.7999
$ws.Range("N136").Value2 = -26272.125
